{"js": "// The document contains a series of \"<Name> is [tab]\" practice-sentence\n// paragraphs (e.g. \"Evan is\", \"Dad is\", \"Daughter is\", \"Usel is\", ...).\n// The edit removes the duplicate/unwanted \"Daughter is\" paragraph entirely\n// (it was an extra line that shouldn't be in the worksheet), while leaving\n// every other paragraph (including \"Dad is\" right before it and \"Usel is\"\n// right after it) completely untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load('items/text');\nawait context.sync();\n\n// Find the paragraph whose text is exactly \"Daughter is\" followed by the\n// paragraph-mark's trailing tab character(s) (Word keeps the tab as part of\n// the paragraph's plain-text representation).\nconst target = paragraphs.items.find((p) => p.text.replace(/\\t+$/, '') === 'Daughter is');\n\nif (target) {\n  target.delete();\n  await context.sync();\n}\n", "ps1": "# The worksheet lists a series of \"<Name> is [tab]\" practice-sentence\n# paragraphs (Evan, Dad, Daughter, Usel, Grandma Sue, Baby Ben, Cindy).\n# This edit removes the extra/duplicate \"Daughter is\" paragraph entirely,\n# leaving every other paragraph (in particular \"Dad is\" right before it and\n# \"Usel is\" right after it) completely untouched.\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"Daughter is\"\n$rng.Find.Forward = $true\n$rng.Find.Wrap = 1\n$found = $rng.Find.Execute()\n\nif ($found) {\n    # Grow the found range to the paragraph that contains it so the tab\n    # character and paragraph mark are removed along with the text,\n    # deleting the paragraph outright instead of leaving it blank.\n    $para = $rng.Paragraphs(1)\n    $para.Range.Delete()\n}\n"}
